$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -207

$ws.Range("H39").Value = 96.333336
$ws.Range("I39").Value = 94.5
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 283.5
$ws.Range("L39").Value = 300
$ws.Range("M39").Value = 12.5
$ws.Range("N39").Value = -892

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws.Range("H132").Value = 5495
$ws.Range("I132").Value = 5495
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16485
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13955

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H45").Value = 1950
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1900
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1523

$ws.Range("H88").Value = 3931.6667
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3931.6667
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3931.6667
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4743.6667

$ws.Range("H91").Value = 3931.6667
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3931.6667
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3931.6667
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6739.6667

$ws.Range("H92").Value = 54750
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 54750
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 54750
$ws.Range("N92").Value = -59742

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 451.75
$ws.Range("I29").Value = 451.75
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 451.75
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -162.75
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 642.63635
$ws.Range("I22").Value = 642.63635
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 642.63635
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -292.63635
$ws.Range("N22").ClearContents()

$ws.Range("H58").Value = 2285.5715
$ws.Range("I58").Value = 2300
$ws.Range("J58").Value = 2249.5
$ws.Range("K58").Value = 2300
$ws.Range("L58").Value = 2249.5
$ws.Range("M58").Value = -2097

$ws.Range("H86").Value = 8025.3335
$ws.Range("I86").Value = 8323.888999999999
$ws.Range("J86").Value = 7129.6665
$ws.Range("K86").Value = 8323.888999999999
$ws.Range("L86").Value = 7129.6665
$ws.Range("M86").Value = -7200.888999999999

$ws.Range("H89").Value = 8025.3335
$ws.Range("I89").Value = 8323.888999999999
$ws.Range("J89").Value = 7129.6665
$ws.Range("K89").Value = 41619.44499999999
$ws.Range("L89").Value = 35648.3325
$ws.Range("M89").Value = -36003.44499999999

$ws.Range("H99").Value = 600799.8
$ws.Range("I99").Value = 999.6667
$ws.Range("J99").Value = 1500500
$ws.Range("K99").Value = 999.6667
$ws.Range("L99").Value = 1500500
$ws.Range("M99").Value = 498.3333
$ws.Range("N99").Value = -1503496

$ws.Range("H114").Value = 58997
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 58997
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 58997
$ws.Range("N114").Value = -67675

$ws.Range("H126").Value = 600799.8
$ws.Range("I126").Value = 999.6667
$ws.Range("J126").Value = 1500500
$ws.Range("K126").Value = 2999.0001
$ws.Range("L126").Value = 4501500
$ws.Range("M126").Value = -529.0001000000002
$ws.Range("N126").Value = -4506440

$ws.Range("H134").Value = 1837.6364
$ws.Range("I134").Value = 1792.375
$ws.Range("J134").Value = 1958.3334
$ws.Range("K134").Value = 5377.125
$ws.Range("L134").Value = 5875.0002
$ws.Range("M134").Value = -2842.125

$ws.Range("H136").Value = 2285.5715
$ws.Range("I136").Value = 2300
$ws.Range("J136").Value = 2249.5
$ws.Range("K136").Value = 6900
$ws.Range("L136").Value = 6748.5
$ws.Range("M136").Value = -4350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 88.73077000000001
$ws.Range("I2").Value = 6.8333335
$ws.Range("J2").Value = 273
$ws.Range("K2").Value = 41.000001
$ws.Range("L2").Value = 1638
$ws.Range("M2").Value = 71.999999

$ws.Range("H12").Value = 178.77777
$ws.Range("I12").Value = 35
$ws.Range("J12").Value = 219.85715
$ws.Range("K12").Value = 105
$ws.Range("L12").Value = 659.5714499999999
$ws.Range("M12").Value = 68
$ws.Range("N12").Value = -1005.57145

$ws.Range("H39").Value = 19166.666
$ws.Range("I39").Value = 3500
$ws.Range("J39").Value = 27000
$ws.Range("K39").Value = 10500
$ws.Range("L39").Value = 81000
$ws.Range("M39").Value = -10206

$ws.Range("H40").Value = 418.63635
$ws.Range("I40").Value = 423.33334
$ws.Range("J40").Value = 397.5
$ws.Range("K40").Value = 1693.33336
$ws.Range("L40").Value = 1590
$ws.Range("M40").Value = -1624.33336
$ws.Range("N40").Value = -1728

$ws.Range("H46").Value = 4862.5
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 6383.3335
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 19150.0005
$ws.Range("M46").Value = -809
$ws.Range("N46").Value = -19332.0005

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 804
$ws.Range("I7").Value = 804
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 804
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -692

$ws.Range("H40").Value = 1730001.6
$ws.Range("I40").Value = 100000
$ws.Range("J40").Value = 2545002.5
$ws.Range("K40").Value = 100000
$ws.Range("L40").Value = 2545002.5
$ws.Range("M40").Value = -99864
$ws.Range("N40").Value = -2545274.5

$ws.Range("H126").Value = 804
$ws.Range("I126").Value = 804
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2412
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 58

$ws.Range("H136").Value = 6133.3335
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -15450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12530
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1296.4
$ws.Range("I132").Value = 1120.5
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3361.5
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -831.5
